# Team_Summary.pptx edit
# Commit message: "Name added to presentation"
#
# 1. The "Challenge" text box (Shapes.Item(3)) still had the placeholder
#    solution name "Sol Name" where the team's actual solution name,
#    "IndiVision", belongs. Replace it in-place, preserving run formatting
#    (bold/italic +mj-lt font) by editing only the matched character range.
# 2. The slide carried an empty, unused <p:timing> animation skeleton
#    (no real effects in MainSequence/InteractiveSequences). Drop it by
#    adding a throwaway effect and immediately deleting it, which prunes
#    the now-empty timing tree from the slide XML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Rename the placeholder solution name -> IndiVision -----------------
$challengeBox = $s.Shapes.Item(3)
$tr = $challengeBox.TextFrame.TextRange
$fullText = $tr.Text

$placeholder = "Sol Name"
$idx = $fullText.IndexOf($placeholder)
if ($idx -ge 0) {
    # Characters() is 1-based.
    $target = $tr.Characters($idx + 1, $placeholder.Length)
    $target.Text = "IndiVision"
}

# --- 2) Remove the empty <p:timing> block -----------------------------------
$timeline = $s.TimeLine
if ($timeline.MainSequence.Count -eq 0 -and $timeline.InteractiveSequences.Count -eq 0) {
    $anchorShape = $s.Shapes.Item(1)
    $dummyEffect = $timeline.MainSequence.AddEffect($anchorShape, 1)
    $timeline.MainSequence.Item(1).Delete()
}
